# Apply the target edits to the workbook (jobs/queries.xlsx).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Remove Sheet2 entirely (table only needed one report sheet) ---
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()
$excel.DisplayAlerts = $true

# --- Header row (text unchanged, just re-asserting) ---
$ws1.Range("F1").Value = "param3"

# --- Row 2: r_comparison_caterpillar ---
$ws1.Range("B2").Value = "r_comparison_caterpillar"
$ws1.Range("C2").Value = "select * from a where id = `$PARAM1"
$ws1.Range("D2").Style = "Normal"
$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "17184,15660,16063,15705,16096,15656,15770,15539,15613,15738"
$ws1.Range("E2").Value = "3,4,5,6,7,8,9,10"

# --- Row 3: r_sub_details_accepted_cases ---
$ws1.Range("B3").Value = "r_sub_details_accepted_cases"
$ws1.Range("C3").Value = "select * from a where id = `$PARAM1 and name = `$PARAM2"
$ws1.Range("D3").Style = "Normal"
$ws1.Range("D3").NumberFormat = "@"
$ws1.Range("D3").Value = "17184,15660,16063,15705,16096,15656,15770,15539,15613,15738"
$ws1.Range("E3").Style = "Normal"
$ws1.Range("E3").Value = "3,4,5,6,7,8,9,10"

# --- Row 4: new row, r_submission_details_first_quadrant ---
$ws1.Rows.Item(4).RowHeight = 17
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "r_submission_details_first_quadrant"
$ws1.Range("C4").Style = "Normal"
$ws1.Range("C4").WrapText = $true
$ws1.Range("C4").Value = "select * from c where id = `$PARAM1 and name = `$PARAM2 and age = `$PARAM3"
$ws1.Range("D4").Style = "Normal"
$ws1.Range("D4").NumberFormat = "@"
$ws1.Range("D4").Value = "17184,15660,16063,15705,16096,15656,15770,15539,15613,15738"
$ws1.Range("E4").Value = "3,4,5,6,7,8,9,10"
$ws1.Range("F4").Value = "99,77,88,99"

# --- Column E width shrinks slightly ---
$ws1.Columns.Item(5).ColumnWidth = 12.833333333333334

# --- Selection cursor moves ---
$ws1.Range("B14").Select()
